# mChat Responses.xlsx update
# "Updates for some responses"
#
# Re-arranges/updates the AI FAQ block (rows 41-48 on the "responses" sheet):
#  - Fixes the "AIaas" action name casing
#  - Splits the combined Advantages/Disadvantages/Benefits button-link cell so each
#    question links to the right set of follow-up topics
#  - Prefixes the disadvantages/benefits answer text with a short lead-in sentence
#  - Fixes two action-name typos (stray space after "action_")
#  - Adds a missing space in the "Process analytics types" link text (row 65)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function ConvertFrom-B64Utf16 {
    param([string]$b64)
    $bytes = [System.Convert]::FromBase64String($b64)
    return [System.Text.Encoding]::Unicode.GetString($bytes)
}

# --- Base64 (UTF-16LE) encoded payloads for strings with embedded tabs/newlines/quotes ---
$B46_B64 = "VABoAGUAcwBlACAAYQByAGUAIAB0AGgAZQAgAEQAaQBzAGEAZAB2AGEAbgB0AGEAZwBlAHMAOgA8AGwAaQA+ACAARQB4AHAAZQBuAHMAaQB2AGUAPAAvAGwAaQA+AAkAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgADwAbABpAD4AUgBlAHEAdQBpAHIAZQBzACAAZABlAGUAcAAgAHQAZQBjAGgAbgBpAGMAYQBsACAAZQB4AHAAZQByAHQAaQBzAGUAIAA8AC8AbABpAD4AIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAIAAgACAAPABsAGkAPgBMAGkAbQBpAHQAZQBkACAAcwB1AHAAcABsAHkAIABvAGYAIABxAHUAYQBsAGkAZgBpAGUAZAAgAHcAbwByAGsAZQByAHMAIAB0AG8AIABiAHUAaQBsAGQAIABBAEkAIAB0AG8AbwBsAHMAPAAvAGwAaQA+AAoAPABsAGkAPgBPAG4AbAB5ACAAawBuAG8AdwBzACAAdwBoAGEAdAAgAGkAdAAnAHMAIABiAGUAZQBuACAAcwBoAG8AdwBuADwALwBsAGkAPgAKADwAbABpAD4ATABhAGMAawAgAG8AZgAgAGEAYgBpAGwAaQB0AHkAIAB0AG8AIABnAGUAbgBlAHIAYQBsAGkAegBlACAAZgByAG8AbQAgAG8AbgBlACAAdABhAHMAawAgAHQAbwAgAGEAbgBvAHQAaABlAHIALgA8AC8AbABpAD4A"
$B47_B64 = "VABoAGUAcwBlACAAYQByAGUAIAB0AGgAZQAgAEIAZQBuAGUAZgBpAHQAcwA6ADwAbABpAD4AQwBoAGEAdABiAG8AdABzACAAdQBzAGUAIABBAEkAIAB0AG8AIAB1AG4AZABlAHIAcwB0AGEAbgBkACAAYwB1AHMAdABvAG0AZQByACAAcAByAG8AYgBsAGUAbQBzACAAZgBhAHMAdABlAHIAIABhAG4AZAAgAHAAcgBvAHYAaQBkAGUAIABtAG8AcgBlACAAZQBmAGYAaQBjAGkAZQBuAHQAIABhAG4AcwB3AGUAcgBzADwALwBsAGkAPgAKADwAbABpAD4ASQBuAHQAZQBsAGwAaQBnAGUAbgB0ACAAYQBzAHMAaQBzAHQAYQBuAHQAcwAgAHUAcwBlACAAQQBJACAAdABvACAAcABhAHIAcwBlACAAYwByAGkAdABpAGMAYQBsACAAaQBuAGYAbwByAG0AYQB0AGkAbwBuACAAZgByAG8AbQAgAGwAYQByAGcAZQAgAGYAcgBlAGUALQB0AGUAeAB0ACAAZABhAHQAYQBzAGUAdABzACAAdABvACAAaQBtAHAAcgBvAHYAZQAgAHMAYwBoAGUAZAB1AGwAaQBuAGcAPAAvAGwAaQA+AAoAPABsAGkAPgBSAGUAYwBvAG0AbQBlAG4AZABhAHQAaQBvAG4AIABlAG4AZwBpAG4AZQBzACAAYwBhAG4AIABwAHIAbwB2AGkAZABlACAAYQB1AHQAbwBtAGEAdABlAGQAIAByAGUAYwBvAG0AbQBlAG4AZABhAHQAaQBvAG4AcwAgAGYAbwByACAAVABWACAAcwBoAG8AdwBzACAAYgBhAHMAZQBkACAAbwBuACAAdQBzAGUAcgBzACAAdgBpAGUAdwBpAG4AZwAgAGgAYQBiAGkAdABzADwALwBsAGkAPgA="

$B46 = ConvertFrom-B64Utf16 $B46_B64
$B47 = ConvertFrom-B64Utf16 $B47_B64

# --- Row 41: action name casing fix (Aiaas -> AIaas) ---
$ws.Range("A41").Value = "action_utter_AI_AIaas"

# --- Row 43: add the advantages/disadvantages/benefits link cell ---
$ws.Range("C43").Value = "Advantages:/AI_Advantages,Disadvantages:/AI_DisAdvantages,Benefits:/AI_Benefits"

# --- Row 44: trim the link cell down to advantages/benefits only ---
$ws.Range("C44").Value = "Advantages:/AI_Advantages,Benefits:/AI_Benefits"

# --- Row 46: AI disadvantages answer now has a lead-in sentence; the link cell moves to row 43/44 ---
$ws.Range("B46").Value = $B46
$ws.Range("C46").ClearContents()

# --- Row 47: fix the action-name typo and add the lead-in sentence to the benefits answer ---
$ws.Range("A47").Value = "action_utter_AI_Benefits"
$ws.Range("B47").Value = $B47

# --- Row 48: fix the action-name typo (stray space) ---
$ws.Range("A48").Value = "action_utter_AI_Technology"

# --- Row 65: add missing space in "Process analytics types" ---
$ws.Range("C65").Value = "Purpose:/ML_Purpose_In_PA, Usages:/How_ML_Used_In_PA,Process analytics types:/PA_Types"

# --- Update the view: scroll position / selection moved while editing this block ---
try {
    $aw = $excel.ActiveWindow
    $aw.ScrollRow = 37
    $aw.ScrollColumn = 1
} catch {
}
$ws.Range("A41").Select()
